$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1407.4
$ws.Range("I28").Value = 669
$ws.Range("J28").Value = 1899.6666
$ws.Range("K28").Value = 669
$ws.Range("L28").Value = 1899.6666
$ws.Range("M28").Value = -184
$ws.Range("N28").Value = -2869.6666
$ws.Range("H80").Value = 1551.3889
$ws.Range("I80").Value = 1192.4
$ws.Range("J80").Value = 2000.125
$ws.Range("K80").Value = 3577.2
$ws.Range("L80").Value = 6000.375
$ws.Range("M80").Value = -2579.2
$ws.Range("N80").Value = -7996.375
$ws.Range("H83").Value = 1551.3889
$ws.Range("I83").Value = 1192.4
$ws.Range("J83").Value = 2000.125
$ws.Range("K83").Value = 10731.6
$ws.Range("L83").Value = 18001.125
$ws.Range("M83").Value = -5739.6
$ws.Range("N83").Value = -27985.125
$ws.Range("H107").Value = 121.42857
$ws.Range("I107").Value = 91.666664
$ws.Range("K107").Value = 91.666664
$ws.Range("M107").Value = 1828.333336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 336.66666
$ws.Range("I63").Value = 336.66666
$ws.Range("K63").Value = 336.66666
$ws.Range("M63").Value = 349.33334
$ws.Range("H66").Value = 336.66666
$ws.Range("I66").Value = 336.66666
$ws.Range("K66").Value = 1683.3333
$ws.Range("M66").Value = 1748.6667
$ws.Range("H74").Value = 100
$ws.Range("I74").Value = 100
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 100
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 774
$ws.Range("N74").Value = ""
$ws.Range("H76").Value = 59333.332
$ws.Range("J76").Value = 59333.332
$ws.Range("L76").Value = 59333.332
$ws.Range("N76").Value = -60009.332
$ws.Range("H77").Value = 100
$ws.Range("I77").Value = 100
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 3868
$ws.Range("N77").Value = ""
$ws.Range("H79").Value = 59333.332
$ws.Range("J79").Value = 59333.332
$ws.Range("L79").Value = 59333.332
$ws.Range("N79").Value = -61673.332
$ws.Range("H102").Value = 2999.75
$ws.Range("I102").Value = 2333
$ws.Range("K102").Value = 2333
$ws.Range("M102").Value = -711
$ws.Range("H110").Value = 123.5
$ws.Range("I110").Value = 123.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 123.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1921.5
$ws.Range("N110").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""
$ws.Range("H107").Value = 2028.8334
$ws.Range("I107").Value = 2028.8334
$ws.Range("K107").Value = 2028.8334
$ws.Range("M107").Value = -108.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 11005
$ws.Range("I62").Value = 11005
$ws.Range("K62").Value = 11005
$ws.Range("M62").Value = -10381
$ws.Range("H65").Value = 11005
$ws.Range("I65").Value = 11005
$ws.Range("K65").Value = 55025
$ws.Range("M65").Value = -51905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 250010.25
$ws.Range("I6").Value = 13.666667
$ws.Range("K6").Value = 41.000001
$ws.Range("M6").Value = 71.999999
$ws.Range("H7").Value = 87
$ws.Range("I7").Value = 87
$ws.Range("K7").Value = 261
$ws.Range("M7").Value = -149
$ws.Range("H131").Value = 1410
$ws.Range("I131").Value = 728.75
$ws.Range("J131").Value = 1864.1666
$ws.Range("K131").Value = 2186.25
$ws.Range("L131").Value = 5592.4998
$ws.Range("M131").Value = 2853.75
$ws.Range("N131").Value = -15672.4998
$ws.Range("H141").Value = 2991.75
$ws.Range("J141").Value = 3000
$ws.Range("L141").Value = 9000
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1457
$ws.Range("I132").Value = 1663.3334
$ws.Range("K132").Value = 4990.0002
$ws.Range("M132").Value = -2460.0002
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1677.1666
$ws.Range("I7").Value = 1677.1666
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1677.1666
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1565.1666
$ws.Range("N7").Value = ""
$ws.Range("H55").Value = 334
$ws.Range("I55").Value = 400
$ws.Range("K55").Value = 400
$ws.Range("M55").Value = -227
$ws.Range("H122").Value = 2302.4285
$ws.Range("I122").Value = 1223.4
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 3670.2
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1220.2
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 1677.1666
$ws.Range("I126").Value = 1677.1666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5031.4998
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2561.4998
$ws.Range("N126").Value = ""
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = ""
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = ""
$ws.Range("N65").Value = ""
$ws.Range("H100").Value = 901.4
$ws.Range("I100").Value = 723.7778
$ws.Range("K100").Value = 1447.5556
$ws.Range("M100").Value = -906.5555999999999

